# bug 59814: clear evaluation workbook and evaluation sheet caches.
# Add structured-reference formulas referencing the table's "Name" column
# on the "Formulas" sheet (rows 2 and 3), and move the active selection
# on the "Table" sheet to A7.

$wb = $excel.ActiveWorkbook

# First move the selection on the "Table" sheet, so that the final
# active/selected sheet ends up being "Formulas" (matching the target
# workbook where the Formulas tab stays selected).
$wsTable = $wb.Worksheets.Item("Table")
$wsTable.Range("A7").Select()

# Now add the two new formula rows on the "Formulas" sheet.
$wsFormulas = $wb.Worksheets.Item("Formulas")
$wsFormulas.Range("A2").Formula = "=\_Prime.1[Name]"
$wsFormulas.Range("A3").Formula = "=\_Prime.1[Name]"

# Leave the selection on A2 of the "Formulas" sheet, and make it the
# active sheet/tab.
$wsFormulas.Range("A2").Select()
